$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("A64").Copy()
$ws.Range("A65").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A65").Value = 45379
$ws.Range("B65").Value = 8
$ws.Range("C65").Value = "Load pre-configured graphs, performance problems"

$ws.Range("C65").Select()
